$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user story header row
$ws.Range("A8").Value = "Kreiranje novcanika"

# Task names (column B) filled first across new rows
$ws.Range("B9").Value = "Aplikativni sloj - Create account metoda"
$ws.Range("B11").Value = "Bank service interface"
$ws.Range("B12").Value = "Bank service mock"
$ws.Range("B13").Value = "Test za kreiranje"

# Time estimates / actuals filled afterwards
$ws.Range("C11").Value = "15min"
$ws.Range("D11").Value = "15min"
$ws.Range("D12").Value = "20min"
$ws.Range("C9").Value = "45min"
$ws.Range("D9").Value = "1h"
$ws.Range("C13").Value = "30min"
$ws.Range("D13").Value = "1h"
$ws.Range("C12").Value = "15min"

# New frontend task added last
$ws.Range("B10").Value = "Frontend - Stranica za create account"

# Widen column A to fit the new, longer content (no longer auto bestFit)
# (ColumnWidth is offset by 5/6 from the stored OOXML width in this runtime,
# so subtract that to land exactly on the target stored width of 19)
$ws.Columns.Item(1).ColumnWidth = 19 - (5/6)

# Move the active selection to reflect where editing left off
$ws.Range("B10").Select()
